# RSCF_QTR_FIN.xlsx - "Doing Updates for Financials"
#
# A new fiscal quarter (period ending 2018-09-30, serial 43373) is appended
# to the RSCF sheet. This is implemented the way a human analyst would do it
# in Excel: insert a brand-new column D (pushing the existing D:K quarters
# out to E:L), carry the formatting of the column that is now to its right,
# and then key in the new quarter's figures down column D. A handful of
# previously-unknown ("NA") cells in the Research & Development row are also
# now known to be zero and get corrected after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSCF")

# --- 1. Insert a new column at D, shifting old D:K to E:L ------------------
$ws.Columns("D:D").Insert(-4161)   # xlShiftToRight

# Carry over number formats/fonts from the column that used to be D (now E)
# so the brand-new column D looks like the rest of the table.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Populate the new column D with the new quarter's data --------------
# Each pair is (row, value). "NA" mirrors the existing placeholder text used
# elsewhere in this sheet for not-yet-available figures.
$newQuarterData = @(
    @(7, 43373),
    @(8, 400),
    @(9, 200),
    @(10, 200),
    @(12, 0),
    @(13, 0),
    @(14, 0),
    @(15, 0),
    @(17, 400),
    @(18, 0),
    @(20, 0),
    @(21, "NA"),
    @(22, 0),
    @(23, 0),
    @(24, 0),
    @(25, 0),
    @(26, 0),
    @(27, 0),
    @(28, 0),
    @(29, 0),
    @(30, 0),
    @(31, 0),
    @(32, 0),
    @(33, 0),
    @(34, 0),
    @(35, 0),
    @(38, 43373),
    @(41, 200),
    @(42, 0),
    @(43, 200),
    @(44, 100),
    @(45, 0),
    @(46, 500),
    @(47, 0),
    @(48, 0),
    @(49, 100),
    @(50, 0),
    @(51, 0),
    @(52, 0),
    @(53, 0),
    @(54, 600),
    @(57, 100),
    @(58, 0),
    @(59, 0),
    @(60, 100),
    @(61, 0),
    @(62, 0),
    @(63, 0),
    @(64, 0),
    @(65, 0),
    @(66, 100),
    @(68, 0),
    @(69, 0),
    @(70, 0),
    @(71, 0),
    @(72, -20200),
    @(73, 0),
    @(74, 0),
    @(75, 0),
    @(76, 500),
    @(77, 0),
    @(80, 43373),
    @(81, 0),
    @(83, 0),
    @(84, 0),
    @(85, 0),
    @(86, 0),
    @(87, 0),
    @(88, 0),
    @(89, 0),
    @(91, 0),
    @(92, 0),
    @(93, 0),
    @(94, 0),
    @(96, 0),
    @(97, 0),
    @(98, 0),
    @(99, 0),
    @(100, 0),
    @(101, 0),
    @(102, 0)
)

foreach ($pair in $newQuarterData) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Range("D$r").Value = $v
}

# --- 3. A few previously-unknown R&D figures are now known to be zero ------
# (Row 14, "Research Development" - these cells shifted from D/E/F/G/I into
# E/F/G/H/J and were placeholder "NA" text; they are corrected to 0 now that
# the new column has been filled in.)
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
